$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.425155276486919
$ws.Range("D2").Value = 0.005047617820290284
$ws.Range("E2").Value = 1.372457414750983
$ws.Range("F2").Value = 0.6848049257048388
$ws.Range("G2").Value = 0.549643369434321
$ws.Range("H2").Value = 0.592691575247386
$ws.Range("L2").Value = 0.5070981716449978
$ws.Range("M2").Value = 0.4074871229341639
$ws.Range("B3").Value = 1.342101873885497
$ws.Range("D3").Value = 0.005286270337466092
$ws.Range("E3").Value = 1.255036682437691
$ws.Range("F3").Value = 0.666708826534574
$ws.Range("G3").Value = 0.5294764152231153
$ws.Range("H3").Value = 0.5921214951306979
$ws.Range("L3").Value = 0.4487414283029523
$ws.Range("M3").Value = 0.3728698081742579
$ws.Range("B4").Value = 1.291359266997659
$ws.Range("D4").Value = 0.005451867117836073
$ws.Range("E4").Value = 1.182892431255084
$ws.Range("F4").Value = 0.6566252842522715
$ws.Range("G4").Value = 0.5180840663426238
$ws.Range("H4").Value = 0.5925421068014742
$ws.Range("L4").Value = 0.4128332772376382
$ws.Range("M4").Value = 0.3516259766708245
$ws.Range("B5").Value = 1.270744948780134
$ws.Range("D5").Value = 0.005524077895589841
$ws.Range("E5").Value = 1.153483006883476
$ws.Range("F5").Value = 0.6527713726784583
$ws.Range("G5").Value = 0.5136867770201548
$ws.Range("H5").Value = 0.5929056630383798
$ws.Range("L5").Value = 0.3981812692121025
$ws.Range("M5").Value = 0.3429718746852757
$ws.Range("B6").Value = 1.267325811739994
$ws.Range("D6").Value = 0.005536352170437908
$ws.Range("E6").Value = 1.148599043951648
$ws.Range("F6").Value = 0.6521467553335469
$ws.Range("G6").Value = 0.5129713074791766
$ws.Range("H6").Value = 0.5929775892375488
$ws.Range("L6").Value = 0.3957471643813903
$ws.Range("M6").Value = 0.3415350475100638
$ws.Range("B7").Value = 1.291080996858085
$ws.Range("D7").Value = 0.005452821917129924
$ws.Range("E7").Value = 1.182495843614561
$ws.Range("F7").Value = 0.6565722799007006
$ws.Range("G7").Value = 0.5180237752676646
$ws.Range("H7").Value = 0.5925462339977798
$ws.Range("L7").Value = 0.4126357523734896
$ws.Range("M7").Value = 0.3515092524353776
$ws.Range("B8").Value = 1.396466137690652
$ws.Range("D8").Value = 0.005125909292816644
$ws.Range("E8").Value = 1.331981654128668
$ws.Range("F8").Value = 0.6783502094545923
$ws.Range("G8").Value = 0.5424819712780504
$ws.Range("H8").Value = 0.5923340794580838
$ws.Range("L8").Value = 0.4869928047463361
$ws.Range("M8").Value = 0.3955487507813658
$ws.Range("B9").Value = 1.605133936623929
$ws.Range("D9").Value = 0.004639057548956771
$ws.Range("E9").Value = 1.624683008847029
$ws.Range("F9").Value = 0.7293557039160277
$ws.Range("G9").Value = 0.5984766374960202
$ws.Range("H9").Value = 0.5981074497658199
$ws.Range("L9").Value = 0.6321976633349209
$ws.Range("M9").Value = 0.4820019249778795
$ws.Range("B10").Value = 1.759691751742423
$ws.Range("D10").Value = 0.004379551020118555
$ws.Range("E10").Value = 1.839405290055566
$ws.Range("F10").Value = 0.7720954121058838
$ws.Range("G10").Value = 0.6447591029819932
$ws.Range("H10").Value = 0.606227460794571
$ws.Range("L10").Value = 0.7385229001829998
$ws.Range("M10").Value = 0.5455855968609455
$ws.Range("B11").Value = 1.830282894892605
$ws.Range("D11").Value = 0.004283739445661183
$ws.Range("E11").Value = 1.937008505043224
$ws.Range("F11").Value = 0.7927272931386966
$ws.Range("G11").Value = 0.6669845914881023
$ws.Range("H11").Value = 0.6107864911029139
$ws.Range("L11").Value = 0.7868200927677833
$ws.Range("M11").Value = 0.5745289821210662
$ws.Range("B12").Value = 1.85705478576341
$ws.Range("D12").Value = 0.004250734362553743
$ws.Range("E12").Value = 1.973956418570253
$ws.Range("F12").Value = 0.80071467796067
$ws.Range("G12").Value = 0.6755735487702452
$ws.Range("H12").Value = 0.6126391423780433
$ws.Range("L12").Value = 0.8050990022102837
$ws.Range("M12").Value = 0.5854919341195881
$ws.Range("B13").Value = 1.851287175697507
$ws.Range("D13").Value = 0.0042576955947915
$ws.Range("E13").Value = 1.965999594727037
$ws.Range("F13").Value = 0.7989866328956481
$ws.Range("G13").Value = 0.6737160173205154
$ws.Range("H13").Value = 0.6122344966595392
$ws.Range("L13").Value = 0.8011627651404751
$ws.Range("M13").Value = 0.5831307449042953
$ws.Range("B14").Value = 1.832484621454967
$ws.Range("D14").Value = 0.004280957971602106
$ws.Range("E14").Value = 1.940048487088234
$ws.Range("F14").Value = 0.7933809018785922
$ws.Range("G14").Value = 0.6676877258679212
$ws.Range("H14").Value = 0.6109363691518013
$ws.Range("L14").Value = 0.7883241138515302
$ws.Range("M14").Value = 0.5754308561368049
$ws.Range("B15").Value = 1.820972796705632
$ws.Range("D15").Value = 0.004295635934166597
$ws.Range("E15").Value = 1.924151038431404
$ws.Range("F15").Value = 0.7899700645070027
$ws.Range("G15").Value = 0.6640178323865769
$ws.Range("H15").Value = 0.6101577240459619
$ws.Range("L15").Value = 0.7804587422570251
$ws.Range("M15").Value = 0.5707148068114378
$ws.Range("B16").Value = 1.755084128605858
$ws.Range("D16").Value = 0.004386266730980637
$ws.Range("E16").Value = 1.833025048765137
$ws.Range("F16").Value = 0.7707713027531327
$ws.Range("G16").Value = 0.6433305447921498
$ws.Range("H16").Value = 0.6059470836223397
$ws.Range("L16").Value = 0.7353651396362295
$ws.Range("M16").Value = 0.5436944582373258
$ws.Range("B17").Value = 1.714735783273795
$ws.Range("D17").Value = 0.004447617252399283
$ws.Range("E17").Value = 1.777101878944848
$ws.Range("F17").Value = 0.7593005998964628
$ws.Range("G17").Value = 0.6309425621135745
$ws.Range("H17").Value = 0.6035868881689339
$ws.Range("L17").Value = 0.707683539869123
$ws.Range("M17").Value = 0.5271231781785986
$ws.Range("B18").Value = 1.691555023816534
$ws.Range("D18").Value = 0.004484992901488738
$ws.Range("E18").Value = 1.744929373245242
$ws.Range("F18").Value = 0.75281469319998
$ws.Range("G18").Value = 0.6239273231780373
$ws.Range("H18").Value = 0.6023106922633303
$ws.Range("L18").Value = 0.6917551609982411
$ws.Range("M18").Value = 0.5175935889264593
$ws.Range("B19").Value = 1.683710976203429
$ws.Range("D19").Value = 0.004498003978241627
$ws.Range("E19").Value = 1.734035165058884
$ws.Range("F19").Value = 0.750637761268834
$ws.Range("G19").Value = 0.6215708488823566
$ws.Range("H19").Value = 0.60189250651419
$ws.Range("L19").Value = 0.6863609440048322
$ws.Range("M19").Value = 0.5143673382930558
$ws.Range("B20").Value = 1.719028186690537
$ws.Range("D20").Value = 0.004440869673610237
$ws.Range("E20").Value = 1.783055729284285
$ws.Range("F20").Value = 0.7605100859859704
$ws.Range("G20").Value = 0.6322498689491738
$ws.Range("H20").Value = 0.6038297048454808
$ws.Range("L20").Value = 0.7106309824678476
$ws.Range("M20").Value = 0.5288870357602065
$ws.Range("B21").Value = 1.838006286444624
$ws.Range("D21").Value = 0.004274035682950483
$ws.Range("E21").Value = 1.947671302043631
$ws.Range("F21").Value = 0.7950226759256651
$ws.Range("G21").Value = 0.6694536619322946
$ws.Range("H21").Value = 0.6113142197357604
$ws.Range("L21").Value = 0.7920954134761473
$ws.Range("M21").Value = 0.5776924261144387
$ws.Range("B22").Value = 1.916001785570359
$ws.Range("D22").Value = 0.004184136784363801
$ws.Range("E22").Value = 2.055184733304259
$ws.Range("F22").Value = 0.8185976439552007
$ws.Range("G22").Value = 0.6947767842823112
$ws.Range("H22").Value = 0.616942498215252
$ws.Range("L22").Value = 0.8452777548096719
$ws.Range("M22").Value = 0.6096054444680306
$ws.Range("B23").Value = 1.874352472539556
$ws.Range("D23").Value = 0.004230340050785486
$ws.Range("E23").Value = 1.997809888540758
$ws.Range("F23").Value = 0.8059208355263792
$ws.Range("G23").Value = 0.6811677008595041
$ws.Range("H23").Value = 0.6138705616926927
$ws.Range("L23").Value = 0.8168987687423339
$ws.Range("M23").Value = 0.592571404974592
$ws.Range("B24").Value = 1.717087541138312
$ws.Range("D24").Value = 0.004443913704688995
$ws.Range("E24").Value = 1.780364060491621
$ws.Range("F24").Value = 0.7599629388482327
$ws.Range("G24").Value = 0.631658503259132
$ws.Range("H24").Value = 0.6037196761650421
$ws.Range("L24").Value = 0.7092984868317842
$ws.Range("M24").Value = 0.5280896036580316
$ws.Range("B25").Value = 1.548466354693119
$ws.Range("D25").Value = 0.004753828655136516
$ws.Range("E25").Value = 1.545552917616192
$ws.Range("F25").Value = 0.7146476117381724
$ws.Range("G25").Value = 0.5824432944151852
$ws.Range("H25").Value = 0.5958721617152491
$ws.Range("L25").Value = 0.5929797532598684
$ws.Range("M25").Value = 0.4586032086208078
